$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.146.54"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "2.893.62"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'528.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "'130.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "2.892.02"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "'0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("D10").Value = "'6.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("E11").Value = "  -4.41%  "
$ws.Range("E12").Value = "  -4.41%  "
$ws.Range("D13").Value = "'0.0000211"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("D14").Value = "'32.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.23%  "
$ws.Range("D15").Value = "3.394.24"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "60.059.93"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "2.900.56"
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").Value = "'449.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("E22").Value = "  -6.28%  "
$ws.Range("D23").Value = "'6.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").Value = "'76.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("D25").Value = "'11.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'2.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "'7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.63%  "
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "'24.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.59%  "
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("D33").Value = "'2.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").Value = "'5.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "'53.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").Value = "'5.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").Value = "'431.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.12%  "
$ws.Range("D38").Value = "'0.0770"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").Value = "'0.0371"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "2.853.54"
$ws.Range("E40").Value = "  -10.33%  "
$ws.Range("E41").Value = "  -6.39%  "
$ws.Range("D42").Value = "'7.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'0.236"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").Value = "'111.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.10%  "
$ws.Range("D50").Value = "0.0₃0467"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "'1.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.07%  "
